# Update movies' (Chinese) names to their English search names on the
# "movie" worksheet, and touch up the formatting for the two rows that
# were just translated (matching the style already used by previously
# translated rows), plus restore the selection location that was active
# when the author finished editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("movie")
$ws.Activate()

$translations = @{
    "B44" = "Offbeat Cops"
    "B45" = "The Survivor"
    "B46" = "The Deal"
    "B47" = "The Confidence Man JP: Episode of the Hero"
    "B48" = "BAD CITY"
    "B49" = "Burial"
    "B50" = "First Oscar"
    "B51" = "The Contractor"
    "B63" = "The Stranger in Our Bed"
    "B64" = "V2 Escape from Hell"
    "B65" = "Old Henry"
    "B70" = "The Advent Calendar"
    "B71" = "The Breitner Commando"
    "B72" = "Silencio"
    "B73" = " The Night Eats the World"
    "B74" = "Murder in the Lens"
    "B75" = "Boy Missing/Secuestro"
    "B76" = "Painkillers"
    "B77" = "Hunting Season"
    "B78" = "Last Breath"
    "B79" = "The Trace"
}

foreach ($addr in $translations.Keys) {
    $ws.Range($addr).Value = $translations[$addr]
}

# Rows 73/74 are newly-translated here (the rest of the list was already
# in English before this edit), so bring their look in line with the
# other already-translated rows: same font (copy format from B3) and the
# same taller row height used throughout the sheet for those rows.
$ws.Range("B3").Copy()
$ws.Range("B73").PasteSpecial(-4122)
$ws.Range("B74").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(73).RowHeight = 16.8
$ws.Rows.Item(74).RowHeight = 16.8

# Leave the view scrolled down near the bottom of the list with the last
# edited cell selected.
$ws.Range("C61").Select()
$excel.ActiveWindow.ScrollRow = 48
